# Grade update - PBL review
# Adds a "Raw Grade" column (M) computed from a weighted mix of homework/quiz
# scores plus a new "PBL" participation bonus column (L).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column M
$ws.Range("M1").Value = "Raw Grade"

# Column M is a bit wider than the other numeric columns to fit the header.
$ws.Columns.Item(13).ColumnWidth = 12.3

# Per-student PBL bonus (L) and raw-grade roll-up (M) formulas.
$ws.Range("L2").Formula = "=0.15+0.1"
$ws.Range("M2").Formula = "=AVERAGE(D2:F2,D2:I2)*0.25+0.15*G2+0.15*J2+0.15*K2+L2"

$ws.Range("L3").Formula = "=0.2+0.1"
$ws.Range("M3").Formula = "=+0.15*G3+0.15*J3+0.15*K3+L3"

$ws.Range("L4").Formula = "=0.2+0.1"
$ws.Range("M4").Formula = "=+0.15*G4+0.15*J4+0.15*K4+L4"

$ws.Range("L5").Formula = "=0.1+0.2"
$ws.Range("M5").Formula = "=+0.15*G5+0.15*J5+0.15*K5+L5"

$ws.Range("L6").Formula = "=0.17+0.08"
$ws.Range("M6").Formula = "=+0.15*G6+0.15*J6+0.15*K6+L6"

$ws.Range("L7").Formula = "=0.1+0.2"
$ws.Range("M7").Formula = "=+0.15*G7+0.15*J7+0.15*K7+L7"

$ws.Range("L8").Formula = "=0.2+0.1"
$ws.Range("M8").Formula = "=+0.15*G8+0.15*J8+0.15*K8+L8"

$ws.Range("L9").Formula = "=0.2+0.1"
$ws.Range("M9").Formula = "=+0.15*G9+0.15*J9+0.15*K9+L9"

$ws.Range("L10").Formula = "=0.2+0.1"
$ws.Range("M10").Formula = "=+0.15*G10+0.15*J10+0.15*K10+L10"

$ws.Range("L11").Formula = "=0.2+0.1"
$ws.Range("M11").Formula = "=+0.15*G11+0.15*J11+0.15*K11+L11"

$ws.Range("L12").Formula = "=0.17+0.08"
$ws.Range("M12").Formula = "=+0.15*G12+0.15*J12+0.15*K12+L12"

$ws.Range("L13").Formula = "=0.2+0.1"
$ws.Range("M13").Formula = "=+0.15*G13+0.15*J13+0.15*K13+L13"

$ws.Range("L14").Formula = "=0.2+0.1"
$ws.Range("M14").Formula = "=+0.15*G14+0.15*J14+0.15*K14+L14"

$ws.Range("L15").Formula = "=0.2+0.1"
$ws.Range("M15").Formula = "=+0.15*G15+0.15*J15+0.15*K15+L15"

# Leave the view focused near the new column, matching where the editor was
# working.
[void]$ws.Range("K2").Select()
